$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Ajo" (Chino / Primera) was added for the
# Macroferia Regional de Talca market. Insert a new row at row 423 so the
# existing historical rows (424-553) shift down by one (to 425-554),
# matching the natural effect of Excel's Rows.Insert, then populate the
# newly inserted row with the new record's values.

$ws.Rows.Item(423).Insert()

$ws.Cells.Item(423, 1).Value = 5
$ws.Cells.Item(423, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(423, 3).Value = "Maule"
$ws.Cells.Item(423, 4).Value = 45215
$ws.Cells.Item(423, 5).Value = 7
$ws.Cells.Item(423, 6).Value = 100112003
$ws.Cells.Item(423, 7).Value = "Ajo"
$ws.Cells.Item(423, 8).Value = "Chino"
$ws.Cells.Item(423, 9).Value = "Primera"
$ws.Cells.Item(423, 10).Value = 300
$ws.Cells.Item(423, 11).Value = 20000
$ws.Cells.Item(423, 12).Value = 20000
$ws.Cells.Item(423, 13).Value = 20000
$ws.Cells.Item(423, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(423, 15).Value = "China"
$ws.Cells.Item(423, 16).Value = 2000
$ws.Cells.Item(423, 17).Value = 10
$ws.Cells.Item(423, 18).Value = "Hortaliza"
